# repull data, push all data, mean calculation
# Update column F ("dSF") values for rows 2-35 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -1
    3  = -8
    4  = -2
    5  = -3
    6  = 0
    7  = -4
    8  = 1
    9  = -8
    10 = 3
    11 = 4
    12 = -3
    13 = -5
    14 = 5
    15 = -4
    16 = -2
    17 = -2
    18 = 0
    19 = 1
    20 = 0
    21 = -2
    22 = 6
    23 = 2
    24 = 0
    25 = 13
    26 = -1
    27 = -4
    28 = 2
    29 = 0
    30 = -1
    31 = -7
    32 = 1
    33 = -2
    34 = 1
    35 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
